$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - rows 2-21
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12430
$ws1.Range("F3").Value = 589
$ws1.Range("F4").Value = 2021
$ws1.Range("F5").Value = 260
$ws1.Range("F8").Value = 12367
$ws1.Range("F9").Value = 3017
$ws1.Range("F13").Value = 17
$ws1.Range("F14").Value = 125
$ws1.Range("F16").Value = 2824
$ws1.Range("F20").Value = 213

# Sheet "全部类型" (sheet4) - rows 2-23
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12430
$ws4.Range("F3").Value = 589
$ws4.Range("F4").Value = 2021
$ws4.Range("F5").Value = 260
$ws4.Range("F9").Value = 12367
$ws4.Range("F10").Value = 3017
$ws4.Range("F14").Value = 17
$ws4.Range("F15").Value = 125
$ws4.Range("F17").Value = 2824
$ws4.Range("F22").Value = 213
